$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 357.5
$ws.Cells.Item(2, 9).Value = 240.625
$ws.Cells.Item(2, 10).Value = 825
$ws.Cells.Item(2, 11).Value = 240.625
$ws.Cells.Item(2, 12).Value = 825
$ws.Cells.Item(2, 13).Value = -127.625
$ws.Cells.Item(2, 14).Value = -1051
$ws.Cells.Item(12, 8).Value = 266.3125
$ws.Cells.Item(12, 9).Value = 266.14285
$ws.Cells.Item(12, 10).Value = 267.5
$ws.Cells.Item(12, 11).Value = 266.14285
$ws.Cells.Item(12, 12).Value = 267.5
$ws.Cells.Item(12, 13).Value = -96.14285000000001
$ws.Cells.Item(12, 14).Value = -607.5
$ws.Cells.Item(18, 8).Value = 2732
$ws.Cells.Item(18, 9).Value = 1924.4445
$ws.Cells.Item(18, 11).Value = 1924.4445
$ws.Cells.Item(18, 13).Value = -1640.4445
$ws.Cells.Item(28, 8).Value = 36505.285
$ws.Cells.Item(28, 9).Value = 37783.63
$ws.Cells.Item(28, 11).Value = 37783.63
$ws.Cells.Item(28, 13).Value = -37298.63
$ws.Cells.Item(41, 8).Value = 59239.35
$ws.Cells.Item(41, 9).Value = 579.8889
$ws.Cells.Item(41, 10).Value = 125231.25
$ws.Cells.Item(41, 11).Value = 579.8889
$ws.Cells.Item(41, 12).Value = 125231.25
$ws.Cells.Item(41, 13).Value = -139.8889
$ws.Cells.Item(41, 14).Value = -126111.25
$ws.Cells.Item(43, 8).Value = 899.75
$ws.Cells.Item(43, 10).Value = 799.5
$ws.Cells.Item(43, 12).Value = 799.5
$ws.Cells.Item(43, 14).Value = -937.5
$ws.Cells.Item(86, 8).Value = 7735504
$ws.Cells.Item(86, 9).Value = 3416
$ws.Cells.Item(86, 11).Value = 3416
$ws.Cells.Item(86, 13).Value = -2293
$ws.Cells.Item(89, 8).Value = 7735504
$ws.Cells.Item(89, 9).Value = 3416
$ws.Cells.Item(89, 11).Value = 17080
$ws.Cells.Item(89, 13).Value = -11464
$ws.Cells.Item(100, 8).Value = 5223.8237
$ws.Cells.Item(100, 9).Value = 6270.5
$ws.Cells.Item(100, 11).Value = 6270.5
$ws.Cells.Item(100, 13).Value = -5729.5
$ws.Cells.Item(125, 8).Value = 2093.7778
$ws.Cells.Item(125, 10).Value = 2124.6667
$ws.Cells.Item(125, 12).Value = 19122.0003
$ws.Cells.Item(125, 14).Value = -24042.0003
$ws.Cells.Item(138, 8).Value = 2146.04
$ws.Cells.Item(138, 9).Value = 1883.9524
$ws.Cells.Item(138, 10).Value = 2335.8276
$ws.Cells.Item(138, 11).Value = 5651.857199999999
$ws.Cells.Item(138, 12).Value = 7007.4828
$ws.Cells.Item(138, 13).Value = -511.8571999999995
$ws.Cells.Item(138, 14).Value = -17287.4828

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 111115420
$ws.Cells.Item(74, 9).Value = 200002780
$ws.Cells.Item(74, 11).Value = 200002780
$ws.Cells.Item(74, 13).Value = -200001906
$ws.Cells.Item(77, 8).Value = 111115420
$ws.Cells.Item(77, 9).Value = 200002780
$ws.Cells.Item(77, 11).Value = 1000013900
$ws.Cells.Item(77, 13).Value = -1000009532
$ws.Cells.Item(132, 8).Value = 3045.4
$ws.Cells.Item(132, 9).Value = 2073.5881
$ws.Cells.Item(132, 11).Value = 6220.7643
$ws.Cells.Item(132, 13).Value = -3690.7643

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(70, 8).Value = 18000
$ws.Cells.Item(70, 10).Value = 18000
$ws.Cells.Item(70, 12).Value = 18000
$ws.Cells.Item(70, 14).Value = -18630
$ws.Cells.Item(73, 8).Value = 18000
$ws.Cells.Item(73, 10).Value = 18000
$ws.Cells.Item(73, 12).Value = 18000
$ws.Cells.Item(73, 14).Value = -20184
$ws.Cells.Item(99, 8).Value = 6899.7144
$ws.Cells.Item(99, 10).Value = 1500
$ws.Cells.Item(99, 12).Value = 1500
$ws.Cells.Item(99, 14).Value = -4496
$ws.Cells.Item(107, 8).Value = 481.2973
$ws.Cells.Item(107, 9).Value = 478.48386
$ws.Cells.Item(107, 10).Value = 495.83334
$ws.Cells.Item(107, 11).Value = 478.48386
$ws.Cells.Item(107, 12).Value = 495.83334
$ws.Cells.Item(107, 13).Value = 1441.51614
$ws.Cells.Item(107, 14).Value = -4335.83334
$ws.Cells.Item(122, 8).Value = 3230.3076
$ws.Cells.Item(122, 9).Value = 3499.889
$ws.Cells.Item(122, 10).Value = 2623.75
$ws.Cells.Item(122, 11).Value = 10499.667
$ws.Cells.Item(122, 12).Value = 7871.25
$ws.Cells.Item(122, 13).Value = -8049.667000000001
$ws.Cells.Item(122, 14).Value = -12771.25
$ws.Cells.Item(126, 8).Value = 6899.7144
$ws.Cells.Item(126, 10).Value = 1500
$ws.Cells.Item(126, 12).Value = 4500
$ws.Cells.Item(126, 14).Value = -9440
$ws.Cells.Item(132, 8).Value = 347149.28
$ws.Cells.Item(132, 9).Value = 1319.3914
$ws.Cells.Item(132, 10).Value = 1672830.5
$ws.Cells.Item(132, 11).Value = 3958.1742
$ws.Cells.Item(132, 12).Value = 5018491.5
$ws.Cells.Item(132, 13).Value = -1428.1742
$ws.Cells.Item(132, 14).Value = -5023551.5
$ws.Cells.Item(134, 8).Value = 2952.8298
$ws.Cells.Item(134, 9).Value = 2584.7073
$ws.Cells.Item(134, 11).Value = 7754.1219
$ws.Cells.Item(134, 13).Value = -5219.1219

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(61, 8).Value = 83.71429000000001
$ws.Cells.Item(61, 9).Value = 70.40000000000001
$ws.Cells.Item(61, 10).Value = 117
$ws.Cells.Item(61, 11).Value = 211.2
$ws.Cells.Item(61, 12).Value = 351
$ws.Cells.Item(61, 13).Value = 3.799999999999983
$ws.Cells.Item(61, 14).Value = -781
$ws.Cells.Item(113, 8).Value = 1897.0588
$ws.Cells.Item(113, 9).Value = 2899.4
$ws.Cells.Item(113, 11).Value = 8698.200000000001
$ws.Cells.Item(113, 13).Value = -6528.200000000001
$ws.Cells.Item(131, 8).Value = 3034.0532
$ws.Cells.Item(131, 9).Value = 1689.4286
$ws.Cells.Item(131, 10).Value = 3142.2415
$ws.Cells.Item(131, 11).Value = 5068.2858
$ws.Cells.Item(131, 12).Value = 9426.7245
$ws.Cells.Item(131, 13).Value = -28.28579999999965
$ws.Cells.Item(131, 14).Value = -19506.7245
$ws.Cells.Item(133, 8).Value = 6385.5386
$ws.Cells.Item(133, 9).Value = 2918.6667
$ws.Cells.Item(133, 10).Value = 9357.143
$ws.Cells.Item(133, 11).Value = 8756.000100000001
$ws.Cells.Item(133, 12).Value = 28071.429
$ws.Cells.Item(133, 13).Value = -3696.000100000001
$ws.Cells.Item(133, 14).Value = -38191.429

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2383.838
$ws.Cells.Item(122, 9).Value = 2029.4231
$ws.Cells.Item(122, 10).Value = 3221.5454
$ws.Cells.Item(122, 11).Value = 6088.2693
$ws.Cells.Item(122, 12).Value = 9664.636200000001
$ws.Cells.Item(122, 13).Value = -3638.2693
$ws.Cells.Item(122, 14).Value = -14564.6362
$ws.Cells.Item(132, 8).Value = 5166.6
$ws.Cells.Item(132, 9).Value = 4357.3
$ws.Cells.Item(132, 10).Value = 6785.2
$ws.Cells.Item(132, 11).Value = 13071.9
$ws.Cells.Item(132, 12).Value = 20355.6
$ws.Cells.Item(132, 13).Value = -10541.9
$ws.Cells.Item(132, 14).Value = -25415.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(5, 8).Value = 12475.25
$ws.Cells.Item(5, 10).Value = 13498.25
$ws.Cells.Item(5, 12).Value = 13498.25
$ws.Cells.Item(5, 14).Value = -13724.25
$ws.Cells.Item(22, 8).Value = 727.2
$ws.Cells.Item(22, 9).Value = 724.5
$ws.Cells.Item(22, 10).Value = 729
$ws.Cells.Item(22, 11).Value = 724.5
$ws.Cells.Item(22, 12).Value = 729
$ws.Cells.Item(22, 13).Value = -429.5
$ws.Cells.Item(22, 14).Value = -1319
$ws.Cells.Item(27, 8).Value = 727.2
$ws.Cells.Item(27, 9).Value = 724.5
$ws.Cells.Item(27, 10).Value = 729
$ws.Cells.Item(27, 11).Value = 724.5
$ws.Cells.Item(27, 12).Value = 729
$ws.Cells.Item(27, 13).Value = -617.5
$ws.Cells.Item(27, 14).Value = -943
$ws.Cells.Item(82, 8).Value = 1611.5333
$ws.Cells.Item(82, 9).Value = 1485.25
$ws.Cells.Item(82, 10).Value = 1755.8572
$ws.Cells.Item(82, 11).Value = 1485.25
$ws.Cells.Item(82, 12).Value = 1755.8572
$ws.Cells.Item(82, 13).Value = -1124.25
$ws.Cells.Item(82, 14).Value = -2477.8572
$ws.Cells.Item(85, 8).Value = 1611.5333
$ws.Cells.Item(85, 9).Value = 1485.25
$ws.Cells.Item(85, 10).Value = 1755.8572
$ws.Cells.Item(85, 11).Value = 1485.25
$ws.Cells.Item(85, 12).Value = 1755.8572
$ws.Cells.Item(85, 13).Value = -237.25
$ws.Cells.Item(85, 14).Value = -4251.8572
$ws.Cells.Item(97, 8).Value = 22500
$ws.Cells.Item(97, 10).Value = 22500
$ws.Cells.Item(97, 12).Value = 22500
$ws.Cells.Item(97, 14).Value = -24482
$ws.Cells.Item(132, 8).Value = 125003656
$ws.Cells.Item(132, 9).Value = 142860540
$ws.Cells.Item(132, 11).Value = 428581620
$ws.Cells.Item(132, 13).Value = -428579090

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3656793.5
$ws.Cells.Item(62, 10).Value = 8071.4287
$ws.Cells.Item(62, 12).Value = 8071.4287
$ws.Cells.Item(62, 14).Value = -9319.4287
$ws.Cells.Item(65, 8).Value = 3656793.5
$ws.Cells.Item(65, 10).Value = 8071.4287
$ws.Cells.Item(65, 12).Value = 40357.14350000001
$ws.Cells.Item(65, 14).Value = -46597.14350000001
$ws.Cells.Item(81, 8).Value = 9096022
$ws.Cells.Item(81, 10).Value = 16674866
$ws.Cells.Item(81, 12).Value = 33349732
$ws.Cells.Item(81, 14).Value = -33351854
$ws.Cells.Item(84, 8).Value = 9096022
$ws.Cells.Item(84, 10).Value = 16674866
$ws.Cells.Item(84, 12).Value = 166748660
$ws.Cells.Item(84, 14).Value = -166759268
$ws.Cells.Item(94, 8).Value = 13712
$ws.Cells.Item(94, 10).Value = 13712
$ws.Cells.Item(94, 12).Value = 13712
$ws.Cells.Item(94, 14).Value = -15514
$ws.Cells.Item(113, 8).Value = 4559.129
$ws.Cells.Item(113, 9).Value = 5452.8
$ws.Cells.Item(113, 10).Value = 2934.2727
$ws.Cells.Item(113, 11).Value = 16358.4
$ws.Cells.Item(113, 12).Value = 8802.8181
$ws.Cells.Item(113, 13).Value = -14188.4
$ws.Cells.Item(113, 14).Value = -13142.8181
$ws.Cells.Item(126, 8).Value = 1355.5883
$ws.Cells.Item(126, 9).Value = 1221.5625
$ws.Cells.Item(126, 11).Value = 3664.6875
$ws.Cells.Item(126, 13).Value = -1194.6875
$ws.Cells.Item(132, 8).Value = 404781.75
$ws.Cells.Item(132, 9).Value = 650305.75
$ws.Cells.Item(132, 10).Value = 4189.9473
$ws.Cells.Item(132, 11).Value = 1950917.25
$ws.Cells.Item(132, 12).Value = 12569.8419
$ws.Cells.Item(132, 13).Value = -1948387.25
$ws.Cells.Item(132, 14).Value = -17629.8419
